$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply the updated "Hot Stock Top 20" values for rows 2-21, columns A-C
$ws.Range("A2").Value = "孚日股份"
$ws.Range("B2").Value = "永泰能源"

$ws.Range("B3").Value = "孚日股份"

$ws.Range("C4").Value = "航天发展"

$ws.Range("A5").Value = "永泰能源"
$ws.Range("B5").Value = "日出东方"

$ws.Range("A6").Value = "龙洲股份"
$ws.Range("B6").Value = "华夏幸福"
$ws.Range("C6").Value = "长城军工"

$ws.Range("A7").Value = "合富中国"
$ws.Range("B7").Value = "平潭发展"
$ws.Range("C7").Value = "日出东方"

$ws.Range("B8").Value = "海南海药"
$ws.Range("C8").Value = "龙洲股份"

$ws.Range("A9").Value = "航天发展"
$ws.Range("B9").Value = "众生药业"
$ws.Range("C9").Value = "合富中国"

$ws.Range("A10").Value = "人民同泰"
$ws.Range("B10").Value = "龙洲股份"
$ws.Range("C10").Value = "众生药业"

$ws.Range("A11").Value = "中水渔业"
$ws.Range("B11").Value = "安泰集团"
$ws.Range("C11").Value = "安泰集团"

$ws.Range("A12").Value = "长城军工"
$ws.Range("B12").Value = "盈新发展"
$ws.Range("C12").Value = "永泰能源"

$ws.Range("A13").Value = "众生药业"
$ws.Range("B13").Value = "长城军工"
$ws.Range("C13").Value = "海马汽车"

$ws.Range("A14").Value = "安泰集团"
$ws.Range("C14").Value = "三木集团"

$ws.Range("A15").Value = "海马汽车"
$ws.Range("B15").Value = "合富中国"
$ws.Range("C15").Value = "东百集团"

$ws.Range("A16").Value = "盈新发展"
$ws.Range("B16").Value = "大东方"
$ws.Range("C16").Value = "胜利股份"

$ws.Range("A17").Value = "胜利股份"
$ws.Range("B17").Value = "中水渔业"
$ws.Range("C17").Value = "多氟多"

$ws.Range("A18").Value = "三木集团"
$ws.Range("B18").Value = "日上集团"
$ws.Range("C18").Value = "人民同泰"

$ws.Range("A19").Value = "海南海药"
$ws.Range("B19").Value = "胜利股份"
$ws.Range("C19").Value = "三花智控"

$ws.Range("A20").Value = "东百集团"
$ws.Range("B20").Value = "三木集团"
$ws.Range("C20").Value = "中水渔业"

$ws.Range("A21").Value = "中国武夷"
$ws.Range("B21").Value = "东百集团"
$ws.Range("C21").Value = "国晟科技"
